# Clean up audio & image files
# - Swap the "phonological" and "semantic" distractor columns (C and D)
#   and relabel the header row to the new dist1/dist2/dist3 naming.
# - Refresh the view (zoom + selection) and the resulting column widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the distractor headers: C was "phonological", D was "semantic".
# After the edit, C holds the semantic distractor and D the phonological one.
$ws.Range("C1").Value = "dist1_semantic_distractor_de"
$ws.Range("D1").Value = "dist2_phonological_distractor_de"
$ws.Range("E1").Value = "dist3_unrelated_distractor_de"

# Swap the C/D column contents for every data row (2-21).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value2
    $dVal = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 3).Value = $dVal
    $ws.Cells.Item($r, 4).Value = $cVal
}

# Resize columns C/D/E to fit the newly swapped/renamed content.
$ws.Columns.Item(3).ColumnWidth = 23.4986979166667
$ws.Columns.Item(4).ColumnWidth = 26.6666666666667
$ws.Columns.Item(5).ColumnWidth = 24.1666666666667

# Refresh the view state: zoom to 140% and move the active selection to C7.
$excel.ActiveWindow.Zoom = 140
$ws.Range("C7").Select()
